$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI TPM-derived values (ligand/receptor expression + edge weights)
# keyed by row number -> hashtable of column letter -> new value
$updates = @{
    2 = @{ "G"=24.440895; "H"=73.32268500000001; "I"=0.2259632546784623; "J"=0.2259632546784623; "M"=1.349272; "N"=4.047815999999999; "O"=0.006367874350338741; "P"=0.006367874350338742; "Q"=32.97741527844; "R"=296.79673750596; "S"=0.001438905613586041; "T"=0.001438905613586041 };
    3 = @{ "G"=24.440895; "H"=73.32268500000001; "I"=0.2259632546784623; "J"=0.2259632546784623; "O"=0.864099120973151; "P"=0.8640991209731511; "Q"=4474.92428184465; "R"=40274.31853660185; "S"=0.1952546497398915; "T"=0.1952546497398915 };
    4 = @{ "G"=24.440895; "H"=73.32268500000001; "I"=0.2259632546784623; "J"=0.2259632546784623; "O"=0.1295330046765102; "P"=0.1295330046765102; "Q"=670.81469458551; "R"=6037.33225126959; "S"=0.02926969932498471; "T"=0.02926969932498472 };
    5 = @{ "I"=0.6065386452756374; "J"=0.6065386452756373; "M"=1.349272; "N"=4.047815999999999; "O"=0.006367874350338741; "P"=0.006367874350338742; "Q"=88.51915686972799; "R"=796.6724118275519; "S"=0.00386236188173994; "T"=0.00386236188173994 };
    6 = @{ "I"=0.6065386452756374; "J"=0.6065386452756373; "O"=0.864099120973151; "P"=0.8640991209731511; "S"=0.5241095102189242; "T"=0.524109510218924 };
    7 = @{ "I"=0.6065386452756374; "J"=0.6065386452756373; "O"=0.1295330046765102; "P"=0.1295330046765102; "S"=0.07856677317497328; "T"=0.07856677317497328 };
    8 = @{ "I"=0.1674981000459004; "J"=0.1674981000459004; "M"=1.349272; "N"=4.047815999999999; "O"=0.006367874350338741; "P"=0.006367874350338742; "Q"=24.44492318639733; "R"=220.0043086775759; "S"=0.001066606855012762; "T"=0.001066606855012762 };
    9 = @{ "I"=0.1674981000459004; "J"=0.1674981000459004; "O"=0.864099120973151; "P"=0.8640991209731511; "S"=0.1447349610143354; "T"=0.1447349610143355 };
    10 = @{ "I"=0.1674981000459004; "J"=0.1674981000459004; "O"=0.1295330046765102; "P"=0.1295330046765102; "S"=0.02169653217655218; "T"=0.02169653217655219 }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}
